# Simplify short-course entry again - only put in defaults
# The "int_timeperiod_shortcourse_mdr" row (row 14) on the "constants" sheet
# is removed entirely; every row below it shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

$ws.Rows("14").Delete() | Out-Null

$ws.Range("A7").Select() | Out-Null
